# Generate Report for Handoff
# Adds a new handoff record (file 916f7834-c89b-4fda-99f8-3a34ec29f5dc.md) as
# row 3 to the "Overview", "zh-cn" and "de-de" worksheets/tables.

$wb = $excel.ActiveWorkbook

$commitHash = "c38c69db8f2410ef00c14c993f24164adda06e70"
$newFile = "916f7834-c89b-4fda-99f8-3a34ec29f5dc.md"
$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = "e2e\" + $newFile
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = "'"
$wsOverview.Range("D3").Style = "Normal"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").NumberFormat = $dateFmt
$wsOverview.Range("G3").Value = "2016-08-17 16:39:33"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/" + $commitHash + "/e2e/" + $newFile,
    "",
    "",
    "e2e\" + $newFile) | Out-Null
$wsOverview.Range("B3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = $newFile
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "916f7834-c89b-4fda-99f8-3a34ec29f5dc.765e995ff3ab725120aba7928cfa5e332c076a07.zh-cn.xlf"
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("H3").Value = "2016-08-17 16:39:29"
$wsZhCn.Range("I3").Value = "'"
$wsZhCn.Range("I3").Style = "Normal"
$wsZhCn.Range("J3").Value = "'"
$wsZhCn.Range("J3").Style = "Normal"
$wsZhCn.Range("K3").NumberFormat = $dateFmt
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L3").Value = "'"
$wsZhCn.Range("L3").Style = "Normal"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("M3").Style = "Normal"
$wsZhCn.Range("N3").Value = "'"
$wsZhCn.Range("N3").Style = "Normal"
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("O3").Style = "Normal"
$wsZhCn.Range("P3").Value = "'"
$wsZhCn.Range("P3").Style = "Normal"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/" + $commitHash + "/e2e/" + $newFile,
    "",
    "",
    $newFile) | Out-Null
$wsZhCn.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = $newFile
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "916f7834-c89b-4fda-99f8-3a34ec29f5dc.765e995ff3ab725120aba7928cfa5e332c076a07.de-de.xlf"
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("H3").Value = "2016-08-17 16:39:33"
$wsDeDe.Range("I3").Value = "'"
$wsDeDe.Range("I3").Style = "Normal"
$wsDeDe.Range("J3").Value = "'"
$wsDeDe.Range("J3").Style = "Normal"
$wsDeDe.Range("K3").NumberFormat = $dateFmt
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L3").Value = "'"
$wsDeDe.Range("L3").Style = "Normal"
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("M3").Style = "Normal"
$wsDeDe.Range("N3").Value = "'"
$wsDeDe.Range("N3").Style = "Normal"
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("O3").Style = "Normal"
$wsDeDe.Range("P3").Value = "'"
$wsDeDe.Range("P3").Style = "Normal"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/" + $commitHash + "/e2e/" + $newFile,
    "",
    "",
    $newFile) | Out-Null
$wsDeDe.Range("A3").Style = "HyperLink"

Write-Host "Overview table range:" $loOverview.Range.Address()
Write-Host "zh-cn table range:" $loZhCn.Range.Address()
Write-Host "de-de table range:" $loDeDe.Range.Address()
